# Update benchmark_class_counts figure data and formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make the header row's "background" cell (B1) bold to match the updated figure styling.
$ws.Range("B1").Font.Bold = $true

# Updated benchmark counts (columns B = background, C = true positives).
$ws.Range("B2").Value = 253
$ws.Range("C2").Value = 40

$ws.Range("B3").Value = 236

$ws.Range("B4").Value = 242
$ws.Range("C4").Value = 42

$ws.Range("B5").Value = 244

$ws.Range("B6").Value = 264

$ws.Range("B7").Value = 238

$ws.Range("B8").Value = 224

$ws.Range("B9").Value = 216

# Leave the last worked cell selected on E6 as in the saved workbook.
$ws.Range("E6").Select()
